# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column after "Week" (column B)
#  - shorten week labels from "W01" style to "W1" style
#  - populate the new Week_Start_Date column with each week's start date
#  - store is_holiday_week as a boolean value instead of a number

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B (ASIN), shifting
# ASIN..is_holiday_week one column to the right (B..I -> C..J).
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week labels (column A) and week-start dates (new column B) for the 16 data rows.
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    # Leading apostrophe forces the date-looking text to stay as text
    # instead of being auto-converted to a date serial number.
    $ws.Cells.Item($row, 2).Value = "'" + $weekStartDates[$i]
    # is_holiday_week (now column J) should be a boolean FALSE, not 0.
    $ws.Cells.Item($row, 10).Value = $false
}
